$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8901190770986138, 0.6722003372141275, 0.7041564859437843, 3.266953414790838)"
$ws.Range("C2").Value = "NIG(1.7263983069214293, 1.340289590232969, 2.7773969637545175, 6.792028590279779)"
$ws.Range("D2").Value = "NIG(0.9473291892256073, 0.6570264562913458, 1.9042417395327373, 3.0015514266393515)"
$ws.Range("E2").Value = "NIG(2.139844638188663, 1.7939016248203832, 2.7858736266966524, 5.667535372604149)"
